$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "dedicacion" column (D2:D26): values were stored as fractions (0.5, 1)
# and should be updated to whole-number percentages (50, 100).
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = $val * 100
    }
}

# Update the active selection on the sheet to I9 (single cell), matching the
# selection stored in the saved file after the edit.
$ws.Range("I9").Select()
